$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 - new shipment: "5410 i5/16/0"
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "5410 i5/16/0"
$ws.Range("B6").Value = 150
$ws.Range("C6").Formula = "=B6*25000"
$ws.Range("D6").Value = 1
$ws.Range("E6").Formula = "=D6*B6"
$ws.Range("H6").Formula = "=E6*25000"

# I6 already carries the "#,##0" styled xf (s=3) from the template; force text
# entry so the dd.m.yy-looking string isn't reinterpreted as a date serial.
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "5.3.24"
$ws.Range("I6").NumberFormat = "#,##0"

# J6/K6 are brand-new cells in this row (no prior xf). J copies the plain
# text style already used by J5; K is left with the default (no) style.
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "271761403878"
$ws.Range("K6").Value = "Cam"

# ---------------------------------------------------------------------------
# Row 7 - new shipment: "7420 i7/16/0"
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "7420 i7/16/0"
$ws.Range("B7").Value = 230
$ws.Range("C7").Formula = "=B7*25000"
$ws.Range("D7").Value = 1
$ws.Range("E7").Formula = "=D7*B7"
$ws.Range("H7").Formula = "=E7*25000"

$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "2.3.24"
$ws.Range("I7").NumberFormat = "#,##0"

$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "9434608205498711249516"
$ws.Range("K7").Value = "Green"

# ---------------------------------------------------------------------------
# Row 8 - new shipment: Dell Precision 7550
# ---------------------------------------------------------------------------
# A8 gains a border + wrap-text text-format style it didn't have before.
$ws.Range("A8").Borders.LineStyle = 1
$ws.Range("A8").WrapText = $true
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "Dell Precision 7550 15.6`"  i7/32/512"

$ws.Range("B8").Value = 429
$ws.Range("C8").Formula = "=B8*25000"
$ws.Range("D8").Value = 1
$ws.Range("E8").Formula = "=D8*B8"
$ws.Range("H8").Formula = "=E8*25000"

$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "2.3.24"
$ws.Range("I8").NumberFormat = "#,##0"

$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "404784410903"
$ws.Range("K8").Value = "CAM"

# ---------------------------------------------------------------------------
# Sheet cosmetics
# ---------------------------------------------------------------------------
# Column widths nudged slightly wider (closest values this engine's 1/6-char
# snapping grid can reach to the authored widths).
$ws.Columns.Item(1).ColumnWidth = 29.5
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 15.67
$ws.Columns.Item(6).ColumnWidth = 11.83
$ws.Columns.Item(8).ColumnWidth = 14.33
$ws.Columns.Item(9).ColumnWidth = 11.17
$ws.Columns.Item(10).ColumnWidth = 21.5

# Page setup: portrait, custom (label-printer) paper size.
$ws.PageSetup.PaperSize = 256
$ws.PageSetup.Orientation = 1

# Final selection, matching the saved cursor position.
$ws.Range("I22").Select()
